$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '30.341.13'
$cell.ClearFormats()
$ws.Range("E2").Value = '  +0.07%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.935.28'
$cell.ClearFormats()
$ws.Range("E3").Value = '  +0.07%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.ClearFormats()
$ws.Range("E4").Value = '  +0.14%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '0.7749'
$cell.ClearFormats()
$ws.Range("E5").Value = '  +8.33%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '246.21'
$cell.ClearFormats()
$ws.Range("E6").Value = '  -2.06%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.9997'
$cell.ClearFormats()
$ws.Range("E7").Value = '  -0.05%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.3207'
$cell.ClearFormats()
$ws.Range("E8").Value = '  -2.73%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '27.82'
$cell.ClearFormats()
$ws.Range("E9").Value = '  +0.66%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.07060'
$cell.ClearFormats()
$ws.Range("E10").Value = '  -2.64%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.7814'
$cell.ClearFormats()
$ws.Range("E11").Value = '  -2.59%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.08007'
$cell.ClearFormats()
$ws.Range("E12").Value = '  -1.14%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '1.932.52'
$cell.ClearFormats()
$ws.Range("E13").Value = '  -0.01%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '5.363'
$cell.ClearFormats()
$ws.Range("E14").Value = '  -2.02%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '94.75'
$cell.ClearFormats()
$ws.Range("E15").Value = '  +0.00%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '14.53'
$cell.ClearFormats()
$ws.Range("E16").Value = '  -3.21%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '30.337.11'
$cell.ClearFormats()
$ws.Range("E17").Value = '  +0.09%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '255.20'
$cell.ClearFormats()
$ws.Range("E18").Value = '  +1.03%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.000007980'
$cell.ClearFormats()
$ws.Range("E19").Value = '  -2.50%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '5.821'
$cell.ClearFormats()
$ws.Range("E20").Value = '  +0.11%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '2.184.11'
$cell.ClearFormats()
$ws.Range("E21").Value = '  +0.04%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.9987'
$cell.ClearFormats()
$ws.Range("E22").Value = '  -0.14%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.9999'
$cell.ClearFormats()
$ws.Range("E23").Value = '  +0.01%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '6.754'
$cell.ClearFormats()
$ws.Range("E24").Value = '  -2.95%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '9.566'
$cell.ClearFormats()
$ws.Range("E25").Value = '  -1.88%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '163.48'
$cell.ClearFormats()
$ws.Range("E26").Value = '  -0.79%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '0.1356'
$cell.ClearFormats()
$ws.Range("E27").Value = '  +4.20%  '

$ws.Range("E28").Value = '  -1.26%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '2.280'
$cell.ClearFormats()
$ws.Range("E29").Value = '  -3.05%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '1.371'
$cell.ClearFormats()
$ws.Range("E30").Value = '  +1.38%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '1.517'
$cell.ClearFormats()
$ws.Range("E31").Value = '  -1.39%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '4.425'
$cell.ClearFormats()
$ws.Range("E32").Value = '  -0.09%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '4.129'
$cell.ClearFormats()
$ws.Range("E33").Value = '  -1.15%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.05174'
$cell.ClearFormats()
$ws.Range("E34").Value = '  -0.71%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.284'
$cell.ClearFormats()
$ws.Range("E35").Value = '  +1.60%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.7508'
$cell.ClearFormats()
$ws.Range("E36").Value = '  +0.53%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '2.772'
$cell.ClearFormats()
$ws.Range("E37").Value = '  -0.18%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.01962'
$cell.ClearFormats()
$ws.Range("E38").Value = '  -0.34%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '2.801'
$cell.ClearFormats()
$ws.Range("E39").Value = '  -0.22%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '78.78'
$cell.ClearFormats()
$ws.Range("E40").Value = '  -0.07%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '6.445'
$cell.ClearFormats()
$ws.Range("E41").Value = '  +0.33%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.4511'
$cell.ClearFormats()
$ws.Range("E42").Value = '  -0.46%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '1.976'
$cell.ClearFormats()
$ws.Range("E43").Value = '  -2.35%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.ClearFormats()
$ws.Range("E44").Value = '  -0.01%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.8357'
$cell.ClearFormats()
$ws.Range("E45").Value = '  -0.81%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '100.81'
$cell.ClearFormats()
$ws.Range("E46").Value = '  -0.65%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '9.789'
$cell.ClearFormats()
$ws.Range("E47").Value = '  +0.22%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '7.506'
$cell.ClearFormats()
$ws.Range("E48").Value = '  +1.01%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '984.67'
$cell.ClearFormats()
$ws.Range("E49").Value = '  +11.25%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '37.25'
$cell.ClearFormats()
$ws.Range("E50").Value = '  +1.26%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.4164'
$cell.ClearFormats()
$ws.Range("E51").Value = '  -0.29%  '
